$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# Update canonical URL on Metadata sheet
$wsMeta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/StructureDefinition/FrRatioMedication"

# Update Date on Metadata sheet
$wsMeta.Range("B8").Value = "2025-05-05T08:11:38+00:00"

# Update Quantity {...} text referencing FrSimpleQuantityMedication on Elements sheet (K5, K6)
$wsElem.Range("K5").Value = "Quantity {https://hl7.fr/ig/fhir/medication/StructureDefinition/FrSimpleQuantityMedication}
"
$wsElem.Range("K6").Value = "Quantity {https://hl7.fr/ig/fhir/medication/StructureDefinition/FrSimpleQuantityMedication}
"

# Re-editing the wrapped-text cells causes the engine to stamp an explicit
# row height on rows 5/6; AutoFit restores the default (no explicit height)
# so the rows keep matching the original (unset) row height.
$wsElem.Rows.Item(5).AutoFit()
$wsElem.Rows.Item(6).AutoFit()

# Update column K width on Elements sheet
# (target stored width 72.69921875 is not exactly reproducible via the
# pixel-quantized ColumnWidth API; 71.8333333333 lands on the closest
# achievable stored width of 72.66666666666667)
$wsElem.Columns.Item(11).ColumnWidth = 71.8333333333
